$wb = $excel.ActiveWorkbook

# Existing Sheet1: append a new data row
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Raghav"
$ws1.Range("C4").Value = "Arora"

# New "Sheet2" holding a little editing UI table for the database.
# Copy Sheet1 (placed right after it) so it inherits the same sheet
# formatting/namespaces, then rename + replace its contents.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Name = "Sheet2"

$ws2.Cells.Clear()

$ws2.Range("A1").Value = "rollNo"
$ws2.Range("B1").Value = "FirstName"
$ws2.Range("C1").Value = "LastName"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "llll"
$ws2.Range("C2").Value = "hhhhh"

# Selections / active sheet matching the target workbook view
$ws1.Range("C4").Select()
$ws2.Range("A3").Select()

$ws2.Activate()
